# Apply Habitat Quality updates: "Flow- Summer Base Flow" (column H) rating
# changes for a set of reaches, plus the recomputed Habitat Quality Scores
# Sum (column P) / Percent (column Q) for the rows whose totals changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column H simply changes from "At Risk" to "Adequate"
# (P/Q totals for these rows are not part of the reported diff, i.e.
# they were already #NUM!/blank or otherwise unaffected).
$hOnlyRows = @(158, 159, 214, 293, 299, 334, 352, 452, 574, 575, 628, 629)
foreach ($r in $hOnlyRows) {
    $ws.Cells.Item($r, 8).Value = "Adequate"
}

# Row 416: "At Risk" -> "Unacceptable"
$ws.Cells.Item(416, 8).Value = "Unacceptable"

# Row 460: "Unacceptable" -> "Adequate"
$ws.Cells.Item(460, 8).Value = "Adequate"

# Row 637: "Unacceptable" -> "At Risk"
$ws.Cells.Item(637, 8).Value = "At Risk"

# Rows where column H changes to "Adequate" AND the Sum (P) / Percent (Q)
# totals are updated accordingly.
$updatedRows = @{
    323 = @{ P = 27; Q = 0.6 }
    324 = @{ P = 35; Q = 0.777777777777778 }
    437 = @{ P = 37; Q = 0.822222222222222 }
    445 = @{ P = 24; Q = 0.533333333333333 }
    657 = @{ P = 35; Q = 0.777777777777778 }
}

foreach ($r in $updatedRows.Keys) {
    $ws.Cells.Item($r, 8).Value = "Adequate"
    $ws.Cells.Item($r, 16).Value = $updatedRows[$r].P
    $ws.Cells.Item($r, 17).Value = $updatedRows[$r].Q
}
